$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to stay text-typed so numeric-looking values
# ("1.02", "230.25", etc.) are not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "41.826.45"
$ws.Range("E2").Value = "  -4.43%  "
$ws.Range("D3").Value = "2.219.66"
$ws.Range("E3").Value = "  -5.43%  "
$ws.Range("D4").Value = "1.02"
$ws.Range("E4").Value = "  +1.37%  "
$ws.Range("D5").Value = "230.25"
$ws.Range("E5").Value = "  -3.74%  "
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  -6.52%  "
$ws.Range("D7").Value = "67.81"
$ws.Range("E7").Value = "  -7.35%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "0.542"
$ws.Range("E9").Value = "  -9.33%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "59.84"
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "0.0963"
$ws.Range("E11").Value = "  -3.74%  "
$ws.Range("D12").Value = "33.67"
$ws.Range("E12").Value = "  +3.17%  "
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").Value = "  -2.60%  "
$ws.Range("D14").Value = "6.63"
$ws.Range("E14").Value = "  -9.51%  "
$ws.Range("D15").Value = "2.573.86"
$ws.Range("E15").Value = "  -4.52%  "
$ws.Range("D16").Value = "14.60"
$ws.Range("E16").Value = "  -10.85%  "
$ws.Range("D17").Value = "0.849"
$ws.Range("E17").Value = "  -5.64%  "
$ws.Range("D18").Value = "2.254.30"
$ws.Range("E18").Value = "  -4.03%  "
$ws.Range("D19").Value = "41.718.63"
$ws.Range("E19").Value = "  -4.55%  "
$ws.Range("D20").Value = "0.0₃0972"
$ws.Range("E20").Value = "  -3.88%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "72.69"
$ws.Range("E21").Value = "  -5.88%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "6.13"
$ws.Range("E22").Value = "  -8.67%  "
$ws.Range("D23").Value = "231.80"
$ws.Range("E23").Value = "  -9.55%  "
$ws.Range("D24").Value = "0.993"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").Value = "3.68"
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("D26").Value = "1.73"
$ws.Range("E26").Value = "  -13.90%  "
$ws.Range("D27").Value = "2.34"
$ws.Range("E27").Value = "  -5.50%  "
$ws.Range("D28").Value = "9.85"
$ws.Range("E28").Value = "  -7.01%  "
$ws.Range("D29").Value = "2.09"
$ws.Range("E29").Value = "  -7.88%  "
$ws.Range("D30").Value = "164.81"
$ws.Range("E30").Value = "  -6.43%  "
$ws.Range("D31").Value = "20.14"
$ws.Range("E31").Value = "  -10.88%  "
$ws.Range("D32").Value = "0.120"
$ws.Range("E32").Value = "  -7.52%  "
$ws.Range("D33").Value = "0.124"
$ws.Range("E33").Value = "  -9.29%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.0698"
$ws.Range("E34").Value = "  -7.64%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "5.17"
$ws.Range("E35").Value = "  -5.37%  "
$ws.Range("D36").Value = "4.68"
$ws.Range("E36").Value = "  -9.36%  "
$ws.Range("D37").Value = "3.44"
$ws.Range("E37").Value = "  -9.20%  "
$ws.Range("D38").Value = "5.95"
$ws.Range("E38").Value = "  -4.91%  "
$ws.Range("D39").Value = "2.18"
$ws.Range("E39").Value = "  -7.39%  "
$ws.Range("D40").Value = "0.0260"
$ws.Range("E40").Value = "  -6.95%  "
$ws.Range("D41").Value = "19.81"
$ws.Range("E41").Value = "  +5.04%  "
$ws.Range("D42").Value = "63.31"
$ws.Range("E42").Value = "  -8.70%  "
$ws.Range("D43").Value = "4.77"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "8.75"
$ws.Range("E44").Value = "  -3.58%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.100"
$ws.Range("E45").Value = "  -9.49%  "
$ws.Range("B46").Value = "BinanceUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D46").Value = "1.01"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.183"
$ws.Range("E47").Value = "  -10.10%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "2.88"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").Value = "1.17"
$ws.Range("E49").Value = "  -4.84%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "2.30"
$ws.Range("E50").Value = "  -8.92%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.506.23"
$ws.Range("E51").Value = "  -2.47%  "

# Restore default styling on column D (clears the temporary Text format)
$priceRange.Style = "Normal"
